$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add 5 new "Event Summary" rows (15-19) below the existing Action rows ---
# Each new row follows the same visual pattern as the existing "ACTION n" header
# rows (e.g. row 6): column A holds a short category label, columns B:M are
# merged into one wrapped-text cell holding the longer narrative text.

$categories = @(
  "EVENT SUMMARY HEADER",
  "ACTION 1 EVENT SUMMARY",
  "ACTION 2 EVENT SUMMARY",
  "ACTION 3 EVENT SUMMARY",
  "ACTION 4 EVENT SUMMARY"
)

$summaries = @(
  "While you were focused on handling these minor affairs in your Kingdom, the disease that is currently affecting the neighbouring Kingdoms have starting to take effect within your own Kingdom.",
  "Having to allocate military personnel to such matters meant that you were lacking manpower to support the efforts in trying to maintain and contain the disease.",
  "Your decision to ignore this issue had caused it to escalate into riots and fighting between these farmers. This had also further accelerated the spread of the disease within the Kingdom.",
  "By giving money to these farmers to quell the issue, you managed to settle this issue for now. However, it came at a cost of the Kingdom's treasury.",
  "Although it had sparked some concerns with the villagers, this was a good move as you managed to get rid the source of this issue. This also meant that you had more time to focus on the disease spreading within the city."
)

$startRow = 15

# Merge the B:M span and copy the formatting from the existing row 6 "ACTION 1"
# banner row onto each new row before writing values, so cell styles line up
# the same way they do for row 6 (merge first, then paste formats, keeps the
# un-merged trailing cell styles intact instead of Excel re-deriving them).
for ($i = 0; $i -lt 5; $i++) {
  $r = $startRow + $i
  $ws.Range("B${r}:M${r}").Merge()
  $ws.Range("A6:M6").Copy()
  $ws.Range("A${r}:M${r}").PasteSpecial(-4122)
  $ws.Rows($r).RowHeight = 41.4
}
$excel.CutCopyMode = $false

# Write the column-A category labels first, then the column-B narrative text,
# so new shared-string entries are appended in that same grouped order.
for ($i = 0; $i -lt 5; $i++) {
  $r = $startRow + $i
  $ws.Range("A${r}").Value = $categories[$i]
}
for ($i = 0; $i -lt 5; $i++) {
  $r = $startRow + $i
  $ws.Range("B${r}").Value = $summaries[$i]
}
